# Auto-generated script applying 2024-10-18 daily crime-data increment
# to the 'violent-crime-full-year' workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6008
$ws.Range("K3").Value = 6194
$ws.Range("C4").Value = 1851
$ws.Range("K4").Value = 1292
$ws.Range("K6").Value = 6818
$ws.Range("C7").Value = 28396
$ws.Range("K7").Value = 20750

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 416
$ws.Range("K6").Value = 461

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 158
$ws.Range("K3").Value = 164
$ws.Range("K7").Value = 460

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 236
$ws.Range("K3").Value = 329
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 901

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 203
$ws.Range("K6").Value = 205
$ws.Range("K7").Value = 706

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 121
$ws.Range("K7").Value = 488

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 182
$ws.Range("K4").Value = 79
$ws.Range("K5").Value = 55
$ws.Range("K7").Value = 607
$ws.Range("K11").Value = 387
$ws.Range("K12").Value = 37
$ws.Range("K15").Value = 214
$ws.Range("K18").Value = 137
$ws.Range("K19").Value = 601
$ws.Range("K20").Value = 493
$ws.Range("K26").Value = 27
$ws.Range("K29").Value = 1132
$ws.Range("K31").Value = 230
$ws.Range("K33").Value = 901
$ws.Range("K36").Value = 265
$ws.Range("K37").Value = 706
$ws.Range("K42").Value = 771
$ws.Range("K44").Value = 176
$ws.Range("K47").Value = 144
$ws.Range("K48").Value = 261
$ws.Range("K49").Value = 114
$ws.Range("K50").Value = 101
$ws.Range("K51").Value = 266
$ws.Range("C63").Value = 280
$ws.Range("K63").Value = 58
$ws.Range("K64").Value = 129
$ws.Range("K65").Value = 488
$ws.Range("K67").Value = 807
$ws.Range("K72").Value = 101
$ws.Range("K79").Value = 520
$ws.Range("K82").Value = 21
$ws.Range("K83").Value = 460
$ws.Range("K84").Value = 164
$ws.Range("K85").Value = 965
$ws.Range("K87").Value = 37
$ws.Range("K91").Value = 234
$ws.Range("K94").Value = 279
$ws.Range("K96").Value = 217
$ws.Range("K98").Value = 102
$ws.Range("C101").Value = 28396
$ws.Range("K101").Value = 20750

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 72
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 289
$ws.Range("K7").Value = 807

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 64
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 66
$ws.Range("K3").Value = 99
$ws.Range("K6").Value = 218

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 326
$ws.Range("K3").Value = 407
$ws.Range("K7").Value = 1132

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 126
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 185
$ws.Range("K6").Value = 193
$ws.Range("K7").Value = 601

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 213
$ws.Range("K3").Value = 232
$ws.Range("K6").Value = 285
$ws.Range("K7").Value = 771

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 66
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 174
$ws.Range("K3").Value = 169
$ws.Range("K6").Value = 130
$ws.Range("K7").Value = 520

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 163
$ws.Range("K3").Value = 160
$ws.Range("K7").Value = 493

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 37
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 79
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 199
$ws.Range("K3").Value = 201
$ws.Range("K6").Value = 163
$ws.Range("K7").Value = 607

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 75
$ws.Range("K4").Value = 22
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 76
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 214

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K3").Value = 15
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 100
$ws.Range("K7").Value = 387

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K2").Value = 12
$ws.Range("K3").Value = 15
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 266

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 313
$ws.Range("K3").Value = 334
$ws.Range("K7").Value = 965

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 11
$ws.Range("K6").Value = 21

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 37
